$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1228.3846
$ws.Range("I28").Value = 1255.75
$ws.Range("K28").Value = 1255.75
$ws.Range("M28").Value = -770.75
$ws.Range("H100").Value = 1672
$ws.Range("I100").Value = 927.25
$ws.Range("K100").Value = 927.25
$ws.Range("M100").Value = -386.25
$ws.Range("H125").Value = 500000160
$ws.Range("J125").Value = 250000260
$ws.Range("L125").Value = 2250002340
$ws.Range("N125").Value = -2250007260

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1316.1
$ws.Range("I2").Value = 707.625
$ws.Range("K2").Value = 707.625
$ws.Range("M2").Value = -594.625
$ws.Range("H97").Value = 521.3333
$ws.Range("J97").Value = 549.5
$ws.Range("L97").Value = 549.5
$ws.Range("N97").Value = -1541.5
$ws.Range("H116").Value = 1316.1
$ws.Range("I116").Value = 707.625
$ws.Range("K116").Value = 707.625
$ws.Range("M116").Value = 1586.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1316.1
$ws.Range("I3").Value = 707.625
$ws.Range("K3").Value = 707.625
$ws.Range("M3").Value = -593.625
$ws.Range("H94").Value = 11066.667
$ws.Range("I94").Value = 11066.667
$ws.Range("K94").Value = 11066.667
$ws.Range("M94").Value = -10615.667
$ws.Range("H99").Value = 2497.9375
$ws.Range("J99").Value = 2664.6365
$ws.Range("L99").Value = 2664.6365
$ws.Range("N99").Value = -5660.636500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 898.3333
$ws.Range("I107").Value = 947.75
$ws.Range("K107").Value = 947.75
$ws.Range("M107").Value = 972.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3304.8
$ws.Range("I5").Value = 4006
$ws.Range("K5").Value = 12018
$ws.Range("M5").Value = -11906
$ws.Range("H23").Value = 156.27272
$ws.Range("J23").Value = 223.66667
$ws.Range("L23").Value = 671.00001
$ws.Range("N23").Value = -1141.00001
$ws.Range("H69").Value = 1991.5
$ws.Range("I69").Value = 1991.5
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 5974.5
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -5163.5
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 1991.5
$ws.Range("I72").Value = 1991.5
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 17923.5
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -13867.5
$ws.Range("N72").ClearContents()
$ws.Range("H80").Value = 2445
$ws.Range("J80").Value = 2445
$ws.Range("L80").Value = 7335
$ws.Range("N80").Value = -9207
$ws.Range("H83").Value = 2445
$ws.Range("J83").Value = 2445
$ws.Range("L83").Value = 22005
$ws.Range("N83").Value = -31365
$ws.Range("H134").Value = 13668.571
$ws.Range("J134").Value = 18357
$ws.Range("L134").Value = 55071
$ws.Range("N134").Value = -65211
$ws.Range("H135").Value = 3304.8
$ws.Range("I135").Value = 4006
$ws.Range("K135").Value = 36054
$ws.Range("M135").Value = -33519

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 500149
$ws.Range("I7").Value = 500149
$ws.Range("K7").Value = 500149
$ws.Range("M7").Value = -500037
$ws.Range("H8").Value = 500149
$ws.Range("I8").Value = 500149
$ws.Range("K8").Value = 500149
$ws.Range("M8").Value = -500010
$ws.Range("H12").Value = 8000000
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H46").Value = 8000
$ws.Range("I46").Value = 8000
$ws.Range("K46").Value = 8000
$ws.Range("M46").Value = -7844
$ws.Range("H80").Value = 3443.1333
$ws.Range("I80").Value = 2129.3333
$ws.Range("K80").Value = 2129.3333
$ws.Range("M80").Value = -1131.3333
$ws.Range("H83").Value = 3443.1333
$ws.Range("I83").Value = 2129.3333
$ws.Range("K83").Value = 10646.6665
$ws.Range("M83").Value = -5654.666499999999
$ws.Range("H102").Value = 1453.9474
$ws.Range("I102").Value = 1260.2941
$ws.Range("K102").Value = 1260.2941
$ws.Range("M102").Value = 361.7058999999999
$ws.Range("H132").Value = 1720.4
$ws.Range("I132").Value = 1650.75
$ws.Range("K132").Value = 4952.25
$ws.Range("M132").Value = -2422.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6968.857
$ws.Range("I7").Value = 4093
$ws.Range("J7").Value = 8119.2
$ws.Range("K7").Value = 4093
$ws.Range("L7").Value = 8119.2
$ws.Range("M7").Value = -3981
$ws.Range("N7").Value = -8343.2
$ws.Range("H23").Value = 89999
$ws.Range("I23").Value = 89999
$ws.Range("K23").Value = 89999
$ws.Range("M23").Value = -89769
$ws.Range("H40").Value = 4711.5
$ws.Range("I40").Value = 4498
$ws.Range("K40").Value = 4498
$ws.Range("M40").Value = -4362
$ws.Range("H61").Value = 2923.818
$ws.Range("I61").Value = 2348.8572
$ws.Range("K61").Value = 2348.8572
$ws.Range("M61").Value = -2146.8572
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H93").Value = 3999.3333
$ws.Range("I93").Value = 2998
$ws.Range("J93").Value = 4500
$ws.Range("K93").Value = 2998
$ws.Range("L93").Value = 4500
$ws.Range("M93").Value = -1750
$ws.Range("N93").Value = -6996
$ws.Range("H113").Value = 2923.818
$ws.Range("I113").Value = 2348.8572
$ws.Range("K113").Value = 2348.8572
$ws.Range("M113").Value = -178.8571999999999
$ws.Range("H126").Value = 6968.857
$ws.Range("I126").Value = 4093
$ws.Range("J126").Value = 8119.2
$ws.Range("K126").Value = 12279
$ws.Range("L126").Value = 24357.6
$ws.Range("M126").Value = -9809
$ws.Range("N126").Value = -29297.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()
$ws.Range("H96").Value = 1351
$ws.Range("I96").Value = 1351
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1351
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 22
$ws.Range("N96").ClearContents()
$ws.Range("H126").Value = 4997
$ws.Range("I126").Value = 3456.4
$ws.Range("J126").Value = 5852.8887
$ws.Range("K126").Value = 10369.2
$ws.Range("L126").Value = 17558.6661
$ws.Range("M126").Value = -7899.200000000001
$ws.Range("N126").Value = -22498.6661
